$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 13-16 (CD4 naive/effector/memory subsets): prefix labels with "CD4/"
$ws.Range("A13").Value = "CD4/CCR7+CD45RA-"
$ws.Range("A14").Value = "CD4/CCR7-CD45RA+"
$ws.Range("A15").Value = "CD4/CCR7-CD45RA-"
$ws.Range("A16").Value = "CD4/CCR7+CD45RA+"

# New rows 17-20: matching CD8 (cytotoxic) subsets, not mislabeled as "helper"
$ws.Range("A17").Value = "CD8/CCR7+CD45RA-"
$ws.Range("B17").Value = "CD8"
$ws.Range("C17").Value = "central memory cytotoxic Tcells (CCR7+ , CD45RA-)"

$ws.Range("A18").Value = "CD8/CCR7-CD45RA+"
$ws.Range("B18").Value = "CD8"
$ws.Range("C18").Value = "effector cytotoxic Tcells  (CCR7-  CD45RA+)"

$ws.Range("A19").Value = "CD8/CCR7-CD45RA-"
$ws.Range("B19").Value = "CD8"
$ws.Range("C19").Value = "effector memory cytotoxic Tcells (CCR7- , CD45RA-)"

$ws.Range("A20").Value = "CD8/CCR7+CD45RA+"
$ws.Range("B20").Value = "CD8"
$ws.Range("C20").Value = "naive cytotoxic Tcells (CCR7+ , CD45RA+)"

# Cosmetic: widen column A to fit the longer labels (stored width 18 once
# Excel's character->width padding is applied), and move the selection to
# where the user ended up after entering the new rows.
$ws.Columns("A").ColumnWidth = 17.166666666666668
$ws.Range("B23").Select() | Out-Null
